$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide master + all slide layouts: update the cached date field text
#    from "28/3/2013" to "1/4/2013" (Date Placeholder 3 on each).
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = "1/4/2013"

$layoutCount = $master.CustomLayouts.Count
for ($li = 1; $li -le $layoutCount; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $sh = $layout.Shapes.Item($si)
        if ($sh.Name -eq "Date Placeholder 3") {
            $sh.TextFrame.TextRange.Text = "1/4/2013"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1 ("Github Commands" title slide): retitle + resubtitle.
#    Title:    "Native " + "Git" + " Commands"  ->  "Github" + " Commands"
#    Subtitle: "Commands you can use in " + "GitBash"
#              -> "Stuff you can do with the " + "github" + " app"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$titleTr = $s1.Shapes.Item(1).TextFrame.TextRange
# Original text is "Native Git Commands" (7 + 3 + 10 chars).
# Collapse "Native " + "Git" (chars 1-10) down to just "Github".
$titleTr.Characters(1, 10).Text = "Github"

$subTr = $s1.Shapes.Item(2).TextFrame.TextRange
# Original text is "Commands you can use in GitBash".
# "Commands you can use in " = chars 1-24 ; "GitBash" = chars 25-31.
$subTr.Characters(25, 7).Text = "github"
$subTr.Characters(1, 24).Text = "Stuff you can do with the "
$subTr.InsertAfter(" app") | Out-Null

# ---------------------------------------------------------------------------
# 3) Slide 2 ("Place new file into local repo" list slide): wipe all bullet
#    text from the content placeholder, leaving it blank.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$contentTr = $s2.Shapes.Item(2).TextFrame.TextRange
$paraCount = $contentTr.Paragraphs().Count
for ($i = $paraCount; $i -ge 1; $i--) {
    $contentTr.Paragraphs($i, 1).Delete() | Out-Null
}

# ---------------------------------------------------------------------------
# 4) Remove one of the two trailing blank slides (position 5).
# ---------------------------------------------------------------------------
$p.Slides.Item(5).Delete()
